$d = $word.ActiveDocument

$replacements = @(
    @("721÷3=", "442÷5="),
    @("211÷7=", "121÷7="),
    @("467÷3=", "886÷8="),
    @("282÷7=", "278÷5="),
    @("576÷3=", "267÷4="),
    @("412÷5=", "486÷6="),
    @("856÷8=", "748÷9="),
    @("312÷6=", "875÷7="),
    @("270÷9=", "764÷7="),
    @("236÷7=", "701÷6="),
    @("191÷4=", "894÷8="),
    @("547÷3=", "950÷8="),
    @("730÷2=", "607÷6="),
    @("102÷5=", "998÷8="),
    @("151÷9=", "147÷5="),
    @("893÷2=", "375÷2="),
    @("737÷3=", "377÷8="),
    @("200÷3=", "242÷3="),
    @("587÷3=", "779÷7="),
    @("887÷7=", "139÷6="),
    @("817÷5=", "426÷5="),
    @("671÷8=", "349÷3="),
    @("969÷6=", "790÷2="),
    @("443÷4=", "350÷9="),
    @("950÷3=", "995÷4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
